$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 18; $r -le 23; $r++) {
    $ws.Cells.Item($r, 1).Value = "203.107.1.34"
}
